$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-8 (A:E)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0

$ws.Range("A3").Value = 0.6
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0

$ws.Range("A8").Value = 9.2
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0

# Clear values in rows 9-11 but keep the existing formatting/style (s="2")
$ws.Range("A9:E11").ClearContents()

# Add new row 12 with the same style as the data rows, left empty
$ws.Range("A9:E9").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A12:E12").ClearContents()

# Update selection to D21
$ws.Range("D21").Select()
